# Auto-generated edit script applying the Phantom_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across 8 Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 984.41174
$ws.Range("I28").Value = 775.4286
$ws.Range("K28").Value = 775.4286
$ws.Range("M28").Value = -290.4286
# Row 69
$ws.Range("H69").Value = 12996
$ws.Range("J69").Value = 13495
$ws.Range("L69").Value = 40485
$ws.Range("N69").Value = -42233
# Row 72
$ws.Range("H72").Value = 12996
$ws.Range("J72").Value = 13495
$ws.Range("L72").Value = 121455
$ws.Range("N72").Value = -130191
# Row 98
$ws.Range("H98").Value = 1166.5
$ws.Range("I98").Value = 1045.2727
$ws.Range("K98").Value = 1045.2727
$ws.Range("M98").Value = 452.7273
# Row 109
$ws.Range("H109").Value = 54999
$ws.Range("J109").Value = 54999
$ws.Range("L109").Value = 54999
$ws.Range("N109").Value = -57773
# Row 122
$ws.Range("H122").Value = 1166.5
$ws.Range("I122").Value = 1045.2727
$ws.Range("K122").Value = 3135.8181
$ws.Range("M122").Value = -685.8181
# Row 127
$ws.Range("H127").Value = 1001.875
$ws.Range("I127").Value = 938
$ws.Range("J127").Value = 1108.3334
$ws.Range("K127").Value = 2814
$ws.Range("L127").Value = 3325.0002
$ws.Range("M127").Value = 2146
$ws.Range("N127").Value = -13245.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1474.8096
$ws.Range("I61").Value = 1419.5264
$ws.Range("K61").Value = 1419.5264
$ws.Range("M61").Value = -1207.5264
# Row 122
$ws.Range("H122").Value = 2869.6667
$ws.Range("I122").Value = 2869.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8609.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6159.000100000001
$ws.Range("N122").ClearContents()
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 136
$ws.Range("H136").Value = 1474.8096
$ws.Range("I136").Value = 1419.5264
$ws.Range("K136").Value = 4258.5792
$ws.Range("M136").Value = -1708.5792

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5755.778
$ws.Range("I86").Value = 5971.857
$ws.Range("K86").Value = 5971.857
$ws.Range("M86").Value = -4848.857
# Row 89
$ws.Range("H89").Value = 5755.778
$ws.Range("I89").Value = 5971.857
$ws.Range("K89").Value = 29859.285
$ws.Range("M89").Value = -24243.285
# Row 134
$ws.Range("H134").Value = 3532.1667
$ws.Range("I134").Value = 3626
$ws.Range("K134").Value = 10878
$ws.Range("M134").Value = -8343
# Row 135
$ws.Range("H135").Value = 34000
$ws.Range("J135").Value = 34000
$ws.Range("L135").Value = 34000
$ws.Range("N135").Value = -44140

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 11669333
$ws.Range("J6").Value = 3999
$ws.Range("L6").Value = 3999
$ws.Range("N6").Value = -4225
# Row 22
$ws.Range("H22").Value = 25481.25
$ws.Range("I22").Value = 641.3333
$ws.Range("K22").Value = 641.3333
$ws.Range("M22").Value = -291.3333
# Row 23
$ws.Range("H23").Value = 15
$ws.Range("J23").Value = 15
$ws.Range("L23").Value = 15
$ws.Range("N23").Value = -495
# Row 27
$ws.Range("H27").Value = 15
$ws.Range("J27").Value = 15
$ws.Range("L27").Value = 15
$ws.Range("N27").Value = -399
# Row 31
$ws.Range("H31").Value = 1497.4
$ws.Range("I31").Value = 1246.75
$ws.Range("K31").Value = 1246.75
$ws.Range("M31").Value = -951.75
# Row 34
$ws.Range("H34").Value = 1497.4
$ws.Range("I34").Value = 1246.75
$ws.Range("K34").Value = 1246.75
$ws.Range("M34").Value = -1044.75
# Row 86
$ws.Range("H86").Value = 7998
$ws.Range("I86").Value = 7990
$ws.Range("K86").Value = 7990
$ws.Range("M86").Value = -6867
# Row 89
$ws.Range("H89").Value = 7998
$ws.Range("I89").Value = 7990
$ws.Range("K89").Value = 39950
$ws.Range("M89").Value = -34334
# Row 107
$ws.Range("H107").Value = 1577
$ws.Range("I107").Value = 966
$ws.Range("K107").Value = 966
$ws.Range("M107").Value = 954
# Row 122
$ws.Range("H122").Value = 5874.75
$ws.Range("I122").Value = 6499.6665
$ws.Range("K122").Value = 19498.9995
$ws.Range("M122").Value = -17048.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 4500
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 13500
$ws.Range("N64").Value = -14040
# Row 67
$ws.Range("H67").Value = 4500
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 13500
$ws.Range("N67").Value = -15372
# Row 109
$ws.Range("H109").Value = 1652.9166
$ws.Range("I109").Value = 805
$ws.Range("K109").Value = 2415
$ws.Range("M109").Value = -1375
# Row 113
$ws.Range("H113").Value = 1793.5454
$ws.Range("J113").Value = 1692.375
$ws.Range("L113").Value = 5077.125
$ws.Range("N113").Value = -9417.125
# Row 137
$ws.Range("H137").Value = 7229.1665
$ws.Range("J137").Value = 8631.25
$ws.Range("L137").Value = 25893.75
$ws.Range("N137").Value = -36093.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2916.1428
$ws.Range("I80").Value = 2846.3333
$ws.Range("K80").Value = 2846.3333
$ws.Range("M80").Value = -1848.3333
# Row 83
$ws.Range("H83").Value = 2916.1428
$ws.Range("I83").Value = 2846.3333
$ws.Range("K83").Value = 14231.6665
$ws.Range("M83").Value = -9239.666499999999
# Row 95
$ws.Range("H95").Value = 27250
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 102
$ws.Range("H102").Value = 5899.5
$ws.Range("I102").Value = 5899.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5899.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -4277.5
$ws.Range("N102").ClearContents()
# Row 123
$ws.Range("H123").Value = 18421
$ws.Range("J123").Value = 18421
$ws.Range("L123").Value = 18421
$ws.Range("N123").Value = -23321
# Row 132
$ws.Range("H132").Value = 1998.1111
$ws.Range("I132").Value = 1997.875
$ws.Range("K132").Value = 5993.625
$ws.Range("M132").Value = -3463.625

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5663.385
$ws.Range("I7").Value = 5663.385
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5663.385
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5551.385
$ws.Range("N7").ClearContents()
# Row 64
$ws.Range("H64").Value = 32000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 32000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 99
$ws.Range("H99").Value = 16933
$ws.Range("I99").Value = 16933
$ws.Range("K99").Value = 16933
$ws.Range("M99").Value = -13938
# Row 126
$ws.Range("H126").Value = 5663.385
$ws.Range("I126").Value = 5663.385
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16990.155
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -14520.155
$ws.Range("N126").ClearContents()
# Row 139
$ws.Range("H139").Value = 44999.5
$ws.Range("J139").Value = 44999.5
$ws.Range("L139").Value = 44999.5
$ws.Range("N139").Value = -55279.5

$ws = $wb.Worksheets.Item("WVR")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
# Row 81
$ws.Range("H81").Value = 5566.6665
$ws.Range("I81").Value = 3800
$ws.Range("K81").Value = 7600
$ws.Range("M81").Value = -6539
# Row 84
$ws.Range("H84").Value = 5566.6665
$ws.Range("I84").Value = 3800
$ws.Range("K84").Value = 38000
$ws.Range("M84").Value = -32696
# Row 123
$ws.Range("H123").Value = 37500
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 37500
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 37500
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -47300
